# Revert "tamoc-119: use hash index (#5033)"
# Removes the "visibilityStatus" column from the Patient, Allergy, and
# Diagnosis reference-data sheets, and swaps sheet1's sample row for a
# freshly generated one.

$wb = $excel.ActiveWorkbook

# ---- Sheet "Patient": drop column K (visibilityStatus) ----------------
$wsPatient = $wb.Worksheets.Item("Patient")
$wsPatient.Range("K1").EntireColumn.Delete()

# Replace the sample data row with the new generated values.
$wsPatient.Range("A2").Value = "3cf6cc61-e682-41ba-aa96-7a46ebef428e"
$wsPatient.Range("B2").Value = "NKIX135928"
$wsPatient.Range("C2").Value = "Helena"
$wsPatient.Range("E2").Value = "Gabbrielli"
$wsPatient.Range("F2").Value = "Inoue"
$wsPatient.Range("G2").Value = 38106.5
$wsPatient.Range("I2").Value = "female"

# ---- Sheet "Allergy": drop column D (visibilityStatus) -----------------
$wsAllergy = $wb.Worksheets.Item("Allergy")
$wsAllergy.Range("D1").EntireColumn.Delete()

# ---- Sheet "Diagnosis": drop column D (visibilityStatus) ---------------
$wsDiagnosis = $wb.Worksheets.Item("Diagnosis")
$wsDiagnosis.Range("D1").EntireColumn.Delete()
